$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Cells whose runs get merged into a single run during the resave; the
# visible text itself is unchanged, so re-typing it via Find/Replace
# reproduces the consolidated-run shape Word produces.
Replace-Text "Cookies werden immer mit Localstorage gesucht" "Cookies werden immer mit Localstorage gesucht"
Replace-Text ", und bei onclick abgespielt" ", und bei onclick abgespielt"
Replace-Text "Funktioniert perfekt und ohne Probleme, evtl bei zu vielen Videos Probleme" "Funktioniert perfekt und ohne Probleme, evtl bei zu vielen Videos Probleme"
Replace-Text "Man kann unter Videos kommentieren, comments wieder angezeigt" "Man kann unter Videos kommentieren, comments wieder angezeigt"
Replace-Text "Man kann nicht in einer Session von einem Anderen User joinen" "Man kann nicht in einer Session von einem Anderen User joinen"
Replace-Text "Man kann die UserID in localstorage " "Man kann die UserID in localstorage "
Replace-Text "Man kann seine eigene Videos anzeigen lassen" "Man kann seine eigene Videos anzeigen lassen"
Replace-Text "seine Videos welche geliked wurden Sehen und diese wiederum abspielen" "seine Videos welche geliked wurden Sehen und diese wiederum abspielen"

# The actual content edit: fill in the still-open test result.
Replace-Text "--" "Muss noch bearbeitet werden"
